# Apply updated crypto price/volume data (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" '27.787.29'
Set-TextCell "E2" '  +1.82%  '

# Row 3
Set-TextCell "D3" '1.880.26'
Set-TextCell "E3" '  +1.36%  '

# Row 4
Set-TextCell "E4" '  +0.14%  '

# Row 5
Set-TextCell "D5" '333.09'
Set-TextCell "E5" '  +2.97%  '

# Row 6
Set-TextCell "D6" '1.004'
Set-TextCell "E6" '  +0.14%  '

# Row 7
Set-TextCell "D7" '0.4722'
Set-TextCell "E7" '  +4.44%  '

# Row 8
Set-TextCell "D8" '0.3958'
Set-TextCell "E8" '  +2.57%  '

# Row 9
Set-TextCell "D9" '47.91'
Set-TextCell "E9" '  -1.59%  '

# Row 10
Set-TextCell "D10" '0.08063'
Set-TextCell "E10" '  +1.70%  '

# Row 11
Set-TextCell "E11" '  +2.09%  '

# Row 12
Set-TextCell "E12" '  +4.31%  '

# Row 13
Set-TextCell "D13" '1.877.75'
Set-TextCell "E13" '  +0.68%  '

# Row 14
Set-TextCell "D14" '5.984'
Set-TextCell "E14" '  +1.39%  '

# Row 15
Set-TextCell "D15" '7.138'
Set-TextCell "E15" '  +0.15%  '

# Row 16
Set-TextCell "D16" '1.004'
Set-TextCell "E16" '  +0.17%  '

# Row 17
Set-TextCell "B17" 'ShibaInu'
Set-TextCell "C17" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell "D17" '0.00001050'
Set-TextCell "E17" '  +2.46%  '

# Row 18
Set-TextCell "B18" 'Litecoin'
Set-TextCell "C18" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell "D18" '87.30'
Set-TextCell "E18" '  +1.73%  '

# Row 19
Set-TextCell "D19" '0.06658'
Set-TextCell "E19" '  +1.34%  '

# Row 20
Set-TextCell "D20" '17.23'
Set-TextCell "E20" '  +1.14%  '

# Row 21
Set-TextCell "E21" '  +0.00%  '

# Row 22
Set-TextCell "D22" '27.798.76'
Set-TextCell "E22" '  +1.86%  '

# Row 23
Set-TextCell "D23" '5.537'
Set-TextCell "E23" '  +0.57%  '

# Row 24
Set-TextCell "E24" '  +1.34%  '

# Row 25
Set-TextCell "E25" '  +0.87%  '

# Row 26
Set-TextCell "D26" '2.082.32'
Set-TextCell "E26" '  +0.23%  '

# Row 27
Set-TextCell "D27" '159.56'
Set-TextCell "E27" '  +3.86%  '

# Row 28
Set-TextCell "D28" '20.20'
Set-TextCell "E28" '  +1.65%  '

# Row 29
Set-TextCell "D29" '2.109'
Set-TextCell "E29" '  +2.32%  '

# Row 30
Set-TextCell "D30" '5.585'
Set-TextCell "E30" '  +2.40%  '

# Row 31
Set-TextCell "D31" '122.00'
Set-TextCell "E31" '  +0.69%  '

# Row 32
Set-TextCell "D32" '0.9888'
Set-TextCell "E32" '  +6.11%  '

# Row 33
Set-TextCell "D33" '0.09550'
Set-TextCell "E33" '  +2.72%  '

# Row 34
Set-TextCell "D34" '1.450'
Set-TextCell "E34" '  -0.69%  '

# Row 35
Set-TextCell "D35" '3.598'
Set-TextCell "E35" '  +0.34%  '

# Row 36
Set-TextCell "D36" '5.368'
Set-TextCell "E36" '  +1.96%  '

# Row 37
Set-TextCell "D37" '0.06133'
Set-TextCell "E37" '  +2.34%  '

# Row 38
Set-TextCell "D38" '0.02260'
Set-TextCell "E38" '  +1.75%  '

# Row 39
Set-TextCell "D39" '1.233'
Set-TextCell "E39" '  +1.05%  '

# Row 40
Set-TextCell "D40" '8.156'
Set-TextCell "E40" '  +0.84%  '

# Row 41
Set-TextCell "D41" '0.6034'
Set-TextCell "E41" '  +2.24%  '

# Row 42
Set-TextCell "D42" '1.003'
Set-TextCell "E42" '  +0.07%  '

# Row 43
Set-TextCell "D43" '0.1903'
Set-TextCell "E43" '  +1.04%  '

# Row 44
Set-TextCell "E44" '  +1.83%  '

# Row 45
Set-TextCell "E45" '  -1.24%  '

# Row 46
Set-TextCell "D46" '0.5742'
Set-TextCell "E46" '  +1.86%  '

# Row 47
Set-TextCell "E47" '  +1.38%  '

# Row 48
Set-TextCell "E48" '  +1.93%  '

# Row 49
Set-TextCell "D49" '3.385'
Set-TextCell "E49" '  +0.37%  '

# Row 50
Set-TextCell "D50" '0.06925'
Set-TextCell "E50" '  +2.77%  '

# Row 51
Set-TextCell "D51" '114.37'
Set-TextCell "E51" '  +5.31%  '
